$wb = $excel.ActiveWorkbook

# Update both the "展览" and "全部类型" sheets which contain duplicated data tables.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 192
    $ws.Range("F4").Value = 271
    $ws.Range("F5").Value = 4062
}
